$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range('D2').Value = '34.759.93'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '1.826.87'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('E6').Value = '  +0.58%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -4.35%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').Value = '2.092.01'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '1.842.08'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').Value = '34.767.58'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').Value = '0.0₃0784'
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('E27').Value = '  +2.22%  '
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('E29').Value = '  -7.72%  '
$ws.Range('E31').Value = '  -1.26%  '
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('E33').Value = '  -1.60%  '
$ws.Range('E34').Value = '  +2.88%  '
$ws.Range('E35').Value = '  +6.78%  '
$ws.Range('E36').Value = '  +12.01%  '
$ws.Range('E37').Value = '  +1.59%  '
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('E39').Value = '  +5.98%  '
$ws.Range('D40').Value = '1.336.99'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('D48').Value = '2.006.86'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('E51').Value = '  +13.36%  '

# --- Numeric-looking values that must remain TEXT: set via formula then convert to static value ---
$ws.Range('D5').Formula = "=""230.37"""
$ws.Range('D8').Formula = "=""39.44"""
$ws.Range('D9').Formula = "=""0.326"""
$ws.Range('D11').Formula = "=""0.0986"""
$ws.Range('D14').Formula = "=""11.27"""
$ws.Range('D15').Formula = "=""0.669"""
$ws.Range('D16').Formula = "=""4.61"""
$ws.Range('D20').Formula = "=""240.17"""
$ws.Range('D21').Formula = "=""12.09"""
$ws.Range('D25').Formula = "=""171.45"""
$ws.Range('D26').Formula = "=""7.73"""
$ws.Range('D28').Formula = "=""17.31"""
$ws.Range('D29').Formula = "=""1.51"""
$ws.Range('D32').Formula = "=""3.91"""
$ws.Range('D35').Formula = "=""1.23"""
$ws.Range('D36').Formula = "=""1.44"""
$ws.Range('D38').Formula = "=""90.98"""
$ws.Range('D41').Formula = "=""0.0193"""
$ws.Range('D42').Formula = "=""14.48"""
$ws.Range('D44').Formula = "=""2.26"""
$ws.Range('D46').Formula = "=""0.0522"""
$ws.Range('D47').Formula = "=""6.24"""
$ws.Range('D50').Formula = "=""0.0671"""
$ws.Range('D51').Formula = "=""3.22"""

$fixRange = $ws.Range('D5,D8,D9,D11,D14,D15,D16,D20,D21,D25,D26,D28,D29,D32,D35,D36,D38,D41,D42,D44,D46,D47,D50,D51')
$fixRange.Copy()
$fixRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false
$wb.CalculateFullRebuild()
